$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new columns I ("I0") and J ("IF") with headers styled like the existing header row
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Copy the formatting (style) of the existing header cell H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill data rows 2..33: column I = 1 (constant), column J = same value as column H
for ($r = 2; $r -le 33; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
